# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions scheduled price/volume refresh, plus a few rank swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.019.34"
$ws.Range("E2").Value = "  +4.72%  "
$ws.Range("D3").Value = "2.228.35"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'231.12"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'61.22"
$ws.Range("E7").Value = "  -4.24%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.402"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'59.06"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").Value = "'0.0900"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "2.570.43"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "'15.63"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "'22.03"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "'0.802"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "'5.58"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "2.255.14"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "42.066.09"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "'72.09"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "'6.01"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "'250.87"
$ws.Range("E23").Value = "  +8.21%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.38"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").Value = "'9.64"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").Value = "'0.143"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").Value = "'168.57"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "'20.00"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").Value = "'5.01"
$ws.Range("E34").Value = "  +5.75%  "
$ws.Range("D35").Value = "'4.64"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "'0.0634"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "'6.65"
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("D38").Value = "'3.69"
$ws.Range("E38").Value = "  -4.44%  "
$ws.Range("D39").Value = "'2.35"
$ws.Range("E39").Value = "  -4.60%  "
$ws.Range("D40").Value = "'0.000254"
$ws.Range("E40").Value = "  +29.68%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'0.0240"
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("D43").Value = "'4.81"
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("D44").Value = "'8.54"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("D45").Value = "'1.22"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "'0.0969"
$ws.Range("E46").Value = "  +4.50%  "
$ws.Range("D47").Value = "'98.73"
$ws.Range("E47").Value = "  -4.73%  "
$ws.Range("D48").Value = "1.477.24"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.81"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'16.42"
$ws.Range("E50").Value = "  -7.76%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'52.26"
$ws.Range("E51").Value = "  +4.78%  "

# Cells above were written with a leading quote-prefix so the engine keeps
# them as literal text (matching the source inlineStr cells) instead of
# auto-coercing to numbers; restore the default cell style afterwards so
# no stray numeric/quote-prefix formatting lingers on the cell.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D14", "D15", "D16", "D17", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
